$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual cell values (referencing current/original row numbers,
#     before any rows are deleted) ---

# Row 2 - RM 2
$ws.Range("E2").Value = -7.2

# Row 3 - RM 8
$ws.Range("F3").Value = ""

# Row 4 - RM 9
$ws.Range("F4").Value = 17.97

# Row 5 - RM 14
$ws.Range("F5").Value = ""

# Row 6 - RM 21
$ws.Range("E6").Value = ""

# Row 8 - RM 38
$ws.Range("F8").Value = ""

# Row 12 - RM 81
$ws.Range("E12").Value = -5.3

# Row 14 - RM 90
$ws.Range("E14").Value = ""

# Row 20 - RM 134
$ws.Range("E20").Value = -7.2

# Row 21 - RM 135
$ws.Range("E21").Value = -8.699999999999999

# Row 22 - RM 138
$ws.Range("E22").Value = ""

# Row 23 - RM 140
$ws.Range("E23").Value = ""
$ws.Range("F23").Value = 16.48

# Row 29 - SC 101
$ws.Range("F29").Value = ""

# Row 31 - SC 119
$ws.Range("F31").Value = 18.06

# Row 32 - SC 120
$ws.Range("D32").Value = -13.6

# Row 33 - SC 132
$ws.Range("E33").Value = -8.1

# Row 34 - SC 193
$ws.Range("D34").Value = ""

# Row 35 - SC 232
$ws.Range("E35").Value = -10.7

# --- Remove rows that were dropped from the dataset ---
# Delete bottom-up so row numbers referenced above stay valid.
$ws.Rows(28).Delete()   # "SC 92"
$ws.Rows(26).Delete()   # "RM 232"
